$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Clear formulas/values in AN2:AN21 (keep style), and remove AO column entirely
$ws.Range("AN2:AN21").ClearContents()
$ws.Range("AO1:AO21").Value = $null

# Scroll the sheet view so that column S is the left-most visible column
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 19
